$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 88
$ws.Range("A2:Q2").Copy($ws.Range("A88:Q88"))
$ws.Cells.Item(88, 8).Clear()
$ws.Cells.Item(88, 1).Value = 43332.65293115741
$ws.Cells.Item(88, 2).Value = "Double observer distance"
$ws.Cells.Item(88, 3).Value = "Florida"
$ws.Cells.Item(88, 4).Value = 43326.0
$ws.Cells.Item(88, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(88, 6).Value = "Austen"
$ws.Cells.Item(88, 7).Value = "Aislyn"
$ws.Cells.Item(88, 9).Value = 30.0
$ws.Cells.Item(88, 10).Value = 13.0
$ws.Cells.Item(88, 11).Value = 45.2512
$ws.Cells.Item(88, 12).Value = -95.07023
$ws.Cells.Item(88, 13).Value = 95.0
$ws.Cells.Item(88, 14).Value = 0.002442129625706002
$ws.Cells.Item(88, 15).Value = 0.0016898148169275373
$ws.Cells.Item(88, 16).Value = 0.0027314814797136933
$ws.Cells.Item(88, 17).Value = 1.0

# Row 89
$ws.Range("A2:Q2").Copy($ws.Range("A89:Q89"))
$ws.Cells.Item(89, 8).Clear()
$ws.Cells.Item(89, 1).Value = 43332.653517280094
$ws.Cells.Item(89, 2).Value = "Double observer distance"
$ws.Cells.Item(89, 3).Value = "Florida"
$ws.Cells.Item(89, 4).Value = 43326.0
$ws.Cells.Item(89, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(89, 6).Value = "Aislyn"
$ws.Cells.Item(89, 7).Value = "Austen"
$ws.Cells.Item(89, 9).Value = 30.0
$ws.Cells.Item(89, 10).Value = 14.0
$ws.Cells.Item(89, 11).Value = 45.25428
$ws.Cells.Item(89, 12).Value = -95.06863
$ws.Cells.Item(89, 13).Value = 115.0
$ws.Cells.Item(89, 14).Value = 0.001967592594155576
$ws.Cells.Item(89, 15).Value = 0.0014699074090458453
$ws.Cells.Item(89, 16).Value = 0.004074074073287193
$ws.Cells.Item(89, 17).Value = 1.0

# Row 90
$ws.Range("A2:Q2").Copy($ws.Range("A90:Q90"))
$ws.Cells.Item(90, 8).Clear()
$ws.Cells.Item(90, 1).Value = 43332.654121597225
$ws.Cells.Item(90, 2).Value = "Double observer distance"
$ws.Cells.Item(90, 3).Value = "Florida"
$ws.Cells.Item(90, 4).Value = 43326.0
$ws.Cells.Item(90, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(90, 6).Value = "Austen"
$ws.Cells.Item(90, 7).Value = "Aislyn"
$ws.Cells.Item(90, 9).Value = 30.0
$ws.Cells.Item(90, 10).Value = 15.0
$ws.Cells.Item(90, 11).Value = 45.2569
$ws.Cells.Item(90, 12).Value = -95.0664
$ws.Cells.Item(90, 13).Value = 110.0
$ws.Cells.Item(90, 14).Value = 0.0028125000026193447
$ws.Cells.Item(90, 15).Value = 0.0021064814791316167
$ws.Cells.Item(90, 16).Value = 0.0031134259261307307
$ws.Cells.Item(90, 17).Value = 1.0

# Row 91
$ws.Range("A87:Q87").Copy($ws.Range("A91:Q91"))
$ws.Cells.Item(91, 6).Clear()
$ws.Cells.Item(91, 7).Clear()
$ws.Cells.Item(91, 1).Value = 43332.656396886574
$ws.Cells.Item(91, 2).Value = "Quadrat survey"
$ws.Cells.Item(91, 3).Value = "Florida"
$ws.Cells.Item(91, 4).Value = 43326.0
$ws.Cells.Item(91, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(91, 8).Value = "0.5m x 0.5m"
$ws.Cells.Item(91, 9).Value = 30.0
$ws.Cells.Item(91, 10).Value = 15.0
$ws.Cells.Item(91, 11).Value = 45.25689
$ws.Cells.Item(91, 12).Value = -95.06645
$ws.Cells.Item(91, 13).Value = 120.0
$ws.Cells.Item(91, 14).Value = 0.005428240736364387
$ws.Cells.Item(91, 15).Value = 0.0019328703710925765
$ws.Cells.Item(91, 16).Value = 0.005810185182781424
$ws.Cells.Item(91, 17).Value = 1.0

# Row 92
$ws.Range("A87:Q87").Copy($ws.Range("A92:Q92"))
$ws.Cells.Item(92, 6).Clear()
$ws.Cells.Item(92, 7).Clear()
$ws.Cells.Item(92, 1).Value = 43332.65708788195
$ws.Cells.Item(92, 2).Value = "Quadrat survey"
$ws.Cells.Item(92, 3).Value = "Florida"
$ws.Cells.Item(92, 4).Value = 43326.0
$ws.Cells.Item(92, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(92, 8).Value = "0.5m x 0.5m"
$ws.Cells.Item(92, 9).Value = 30.0
$ws.Cells.Item(92, 10).Value = 14.0
$ws.Cells.Item(92, 11).Value = 45.25429
$ws.Cells.Item(92, 12).Value = -95.06862
$ws.Cells.Item(92, 13).Value = 100.0
$ws.Cells.Item(92, 14).Value = 0.003657407403807156
$ws.Cells.Item(92, 15).Value = 0.001527777778392192
$ws.Cells.Item(92, 16).Value = 0.005474537036207039
$ws.Cells.Item(92, 17).Value = 1.0

# Row 93
$ws.Range("A87:Q87").Copy($ws.Range("A93:Q93"))
$ws.Cells.Item(93, 6).Clear()
$ws.Cells.Item(93, 7).Clear()
$ws.Cells.Item(93, 1).Value = 43332.65768975695
$ws.Cells.Item(93, 2).Value = "Quadrat survey"
$ws.Cells.Item(93, 3).Value = "Florida"
$ws.Cells.Item(93, 4).Value = 43326.0
$ws.Cells.Item(93, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(93, 8).Value = "0.5m x 0.5m"
$ws.Cells.Item(93, 9).Value = 30.0
$ws.Cells.Item(93, 10).Value = 13.0
$ws.Cells.Item(93, 11).Value = 45.25121
$ws.Cells.Item(93, 12).Value = -95.07023
$ws.Cells.Item(93, 13).Value = 80.0
$ws.Cells.Item(93, 14).Value = 0.0034143518496421166
$ws.Cells.Item(93, 15).Value = 0.0017129629632108845
$ws.Cells.Item(93, 16).Value = 0.004826388889341615
$ws.Cells.Item(93, 17).Value = 1.0

# Row 94
$ws.Range("A87:Q87").Copy($ws.Range("A94:Q94"))
$ws.Cells.Item(94, 6).Clear()
$ws.Cells.Item(94, 7).Clear()
$ws.Cells.Item(94, 1).Value = 43332.6584206713
$ws.Cells.Item(94, 2).Value = "Quadrat survey"
$ws.Cells.Item(94, 3).Value = "Florida"
$ws.Cells.Item(94, 4).Value = 43326.0
$ws.Cells.Item(94, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(94, 8).Value = "0.5m x 0.5m"
$ws.Cells.Item(94, 9).Value = 30.0
$ws.Cells.Item(94, 10).Value = 12.0
$ws.Cells.Item(94, 11).Value = 45.24731
$ws.Cells.Item(94, 12).Value = -95.06965
$ws.Cells.Item(94, 13).Value = 65.0
$ws.Cells.Item(94, 14).Value = 0.0037037037036498077
$ws.Cells.Item(94, 15).Value = 0.0016898148169275373
$ws.Cells.Item(94, 16).Value = 0.004513888889050577
$ws.Cells.Item(94, 17).Value = 0.75

# Row 95
$ws.Range("A87:Q87").Copy($ws.Range("A95:Q95"))
$ws.Cells.Item(95, 6).Clear()
$ws.Cells.Item(95, 7).Clear()
$ws.Cells.Item(95, 1).Value = 43332.65911978009
$ws.Cells.Item(95, 2).Value = "Quadrat survey"
$ws.Cells.Item(95, 3).Value = "Florida"
$ws.Cells.Item(95, 4).Value = 43326.0
$ws.Cells.Item(95, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(95, 8).Value = "0.5m x 0.5m"
$ws.Cells.Item(95, 9).Value = 30.0
$ws.Cells.Item(95, 10).Value = 11.0
$ws.Cells.Item(95, 11).Value = 45.24428
$ws.Cells.Item(95, 12).Value = -95.06659
$ws.Cells.Item(95, 13).Value = 10.0
$ws.Cells.Item(95, 14).Value = 0.003726851849933155
$ws.Cells.Item(95, 15).Value = 0.0016666666633682325
$ws.Cells.Item(95, 16).Value = 0.00462962962774327
$ws.Cells.Item(95, 17).Value = 1.0

# Row 96
$ws.Range("A87:Q87").Copy($ws.Range("A96:Q96"))
$ws.Cells.Item(96, 6).Clear()
$ws.Cells.Item(96, 7).Clear()
$ws.Cells.Item(96, 1).Value = 43332.65980047453
$ws.Cells.Item(96, 2).Value = "Quadrat survey"
$ws.Cells.Item(96, 3).Value = "Florida"
$ws.Cells.Item(96, 4).Value = 43327.0
$ws.Cells.Item(96, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(96, 8).Value = "0.5m x 0.5m"
$ws.Cells.Item(96, 9).Value = 30.0
$ws.Cells.Item(96, 10).Value = 9.0
$ws.Cells.Item(96, 11).Value = 45.24718
$ws.Cells.Item(96, 12).Value = -95.05644
$ws.Cells.Item(96, 13).Value = 320.0
$ws.Cells.Item(96, 14).Value = 0.00555555555911269
$ws.Cells.Item(96, 15).Value = 0.001527777778392192
$ws.Cells.Item(96, 16).Value = 0.008182870369637385
$ws.Cells.Item(96, 17).Value = 1.25

# Row 97
$ws.Range("A87:Q87").Copy($ws.Range("A97:Q97"))
$ws.Cells.Item(97, 6).Clear()
$ws.Cells.Item(97, 7).Clear()
$ws.Cells.Item(97, 1).Value = 43332.66041
$ws.Cells.Item(97, 2).Value = "Quadrat survey"
$ws.Cells.Item(97, 3).Value = "Florida"
$ws.Cells.Item(97, 4).Value = 43327.0
$ws.Cells.Item(97, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(97, 8).Value = "0.5m x 0.5m"
$ws.Cells.Item(97, 9).Value = 30.0
$ws.Cells.Item(97, 10).Value = 10.0
$ws.Cells.Item(97, 11).Value = 45.24418
$ws.Cells.Item(97, 12).Value = -95.06046
$ws.Cells.Item(97, 13).Value = 335.0
$ws.Cells.Item(97, 14).Value = 0.004583333335176576
$ws.Cells.Item(97, 15).Value = 0.001747685186273884
$ws.Cells.Item(97, 16).Value = 0.00555555555911269
$ws.Cells.Item(97, 17).Value = 1.25

# Row 98
$ws.Range("A2:Q2").Copy($ws.Range("A98:Q98"))
$ws.Cells.Item(98, 8).Clear()
$ws.Cells.Item(98, 1).Value = 43332.66148480324
$ws.Cells.Item(98, 2).Value = "Double observer no distance"
$ws.Cells.Item(98, 3).Value = "Florida"
$ws.Cells.Item(98, 4).Value = 43327.0
$ws.Cells.Item(98, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(98, 6).Value = "Austen"
$ws.Cells.Item(98, 7).Value = "Aislyn"
$ws.Cells.Item(98, 9).Value = 30.0
$ws.Cells.Item(98, 10).Value = 11.0
$ws.Cells.Item(98, 11).Value = 45.24429
$ws.Cells.Item(98, 12).Value = -95.06665
$ws.Cells.Item(98, 13).Value = 25.0
$ws.Cells.Item(98, 14).Value = 0.003969907411374152
$ws.Cells.Item(98, 15).Value = 0.0015625000014551915
$ws.Cells.Item(98, 16).Value = 0.003229166664823424
$ws.Cells.Item(98, 17).Value = 1.0

# Row 99
$ws.Range("A2:Q2").Copy($ws.Range("A99:Q99"))
$ws.Cells.Item(99, 8).Clear()
$ws.Cells.Item(99, 1).Value = 43332.662280046294
$ws.Cells.Item(99, 2).Value = "Double observer no distance"
$ws.Cells.Item(99, 3).Value = "Florida"
$ws.Cells.Item(99, 4).Value = 43327.0
$ws.Cells.Item(99, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(99, 6).Value = "Aislyn"
$ws.Cells.Item(99, 7).Value = "Austen"
$ws.Cells.Item(99, 9).Value = 30.0
$ws.Cells.Item(99, 10).Value = 12.0
$ws.Cells.Item(99, 11).Value = 45.2474
$ws.Cells.Item(99, 12).Value = -95.06975
$ws.Cells.Item(99, 13).Value = 60.0
$ws.Cells.Item(99, 14).Value = 0.003599537034460809
$ws.Cells.Item(99, 15).Value = 0.001631944440305233
$ws.Cells.Item(99, 16).Value = 0.003506944442051463
$ws.Cells.Item(99, 17).Value = 1.25

# Row 100
$ws.Range("A2:Q2").Copy($ws.Range("A100:Q100"))
$ws.Cells.Item(100, 8).Clear()
$ws.Cells.Item(100, 1).Value = 43332.66296215278
$ws.Cells.Item(100, 2).Value = "Double observer no distance"
$ws.Cells.Item(100, 3).Value = "Florida"
$ws.Cells.Item(100, 4).Value = 43327.0
$ws.Cells.Item(100, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(100, 6).Value = "Austen"
$ws.Cells.Item(100, 7).Value = "Aislyn"
$ws.Cells.Item(100, 9).Value = 30.0
$ws.Cells.Item(100, 10).Value = 13.0
$ws.Cells.Item(100, 11).Value = 45.25124
$ws.Cells.Item(100, 12).Value = -95.07021
$ws.Cells.Item(100, 13).Value = 80.0
$ws.Cells.Item(100, 14).Value = 0.0036921296268701553
$ws.Cells.Item(100, 15).Value = 0.0019907407404389232
$ws.Cells.Item(100, 16).Value = 0.003912037034751847
$ws.Cells.Item(100, 17).Value = 1.5

# Row 101
$ws.Range("A2:Q2").Copy($ws.Range("A101:Q101"))
$ws.Cells.Item(101, 8).Clear()
$ws.Cells.Item(101, 1).Value = 43332.66386605324
$ws.Cells.Item(101, 2).Value = "Double observer no distance"
$ws.Cells.Item(101, 3).Value = "Florida"
$ws.Cells.Item(101, 4).Value = 43327.0
$ws.Cells.Item(101, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(101, 6).Value = "Aislyn"
$ws.Cells.Item(101, 7).Value = "Austen"
$ws.Cells.Item(101, 9).Value = 30.0
$ws.Cells.Item(101, 10).Value = 14.0
$ws.Cells.Item(101, 11).Value = 45.25433
$ws.Cells.Item(101, 12).Value = -95.0686
$ws.Cells.Item(101, 13).Value = 90.0
$ws.Cells.Item(101, 14).Value = 0.003564814818673767
$ws.Cells.Item(101, 15).Value = 0.002025462963501923
$ws.Cells.Item(101, 16).Value = 0.0033333333340124227
$ws.Cells.Item(101, 17).Value = 1.25

# Row 102
$ws.Range("A2:Q2").Copy($ws.Range("A102:Q102"))
$ws.Cells.Item(102, 8).Clear()
$ws.Cells.Item(102, 1).Value = 43332.66446453704
$ws.Cells.Item(102, 2).Value = "Double observer no distance"
$ws.Cells.Item(102, 3).Value = "Florida"
$ws.Cells.Item(102, 4).Value = 43327.0
$ws.Cells.Item(102, 5).Value = "Aislyn, Austen"
$ws.Cells.Item(102, 6).Value = "Austen"
$ws.Cells.Item(102, 7).Value = "Aislyn"
$ws.Cells.Item(102, 9).Value = 30.0
$ws.Cells.Item(102, 10).Value = 15.0
$ws.Cells.Item(102, 11).Value = 45.25692
$ws.Cells.Item(102, 12).Value = -95.06639
$ws.Cells.Item(102, 13).Value = 120.0
$ws.Cells.Item(102, 14).Value = 0.003958333334594499
$ws.Cells.Item(102, 15).Value = 0.0017592592630535364
$ws.Cells.Item(102, 16).Value = 0.003125000002910383
$ws.Cells.Item(102, 17).Value = 1.25

